$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cálculos Sectores")

# Fill in the error-balance results for the "by sectors" calculations tab
$ws.Range("C7").Value = "falta"
$ws.Range("C8").Value = "Listo"
$ws.Range("C9").Value = "algo falla"
$ws.Range("C10").Value = "algo falla"
$ws.Range("D10").Value = "Los valores de izquierda son muy altos en el balance…"
$ws.Range("C11").Value = "Trabajando en ello"

# Update page setup (paper size A4 / portrait) for this sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Reflect the scrolled/selected view state recorded in the saved file
$ws.Activate()
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 7
